$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for 2025/11/17.
# Force column A to be stored as text (not auto-parsed as a date serial)
# by pre-setting a text number format before assigning the value.
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2025/11/17"
$ws.Range("B8").Value = "逃离鸭科夫"
$ws.Range("C8").Value = 1155

# Copy the formatting (center/center alignment) used by the other data
# rows onto the new row so it matches the existing style.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)
